$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (the "prueba" and "KNN" rows) entirely
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Remove the now-unused I and J columns entirely
$ws.Columns.Item(10).Delete()
$ws.Columns.Item(9).Delete()

# Update header row
$ws.Range("B1").Value = "precision_macro"
$ws.Range("C1").Value = "recall_macro"
$ws.Range("D1").Value = "precision_micro"
$ws.Range("E1").Value = "recall_micro"
$ws.Range("F1").Value = "f1_macro"
$ws.Range("G1").Value = "accuracy"
$ws.Range("H1").Value = "roc_auc"

# Row 2 - Regresion Logistica RGB
$ws.Range("A2").Value = "Regresión Logística RGB"
$ws.Range("B2").Value = 0.37615196453386
$ws.Range("C2").Value = 0.3703476412684976
$ws.Range("D2").Value = 0.3788135593220339
$ws.Range("E2").Value = 0.3788135593220339
$ws.Range("F2").Value = 0.3689149283290373
$ws.Range("G2").Value = 0.3788135593220339
$ws.Range("H2").Value = ""

# Row 3 - Regresion Logistica RGB
$ws.Range("A3").Value = "Regresión Logística RGB"
$ws.Range("B3").Value = 0.37615196453386
$ws.Range("C3").Value = 0.3703476412684976
$ws.Range("D3").Value = 0.3788135593220339
$ws.Range("E3").Value = 0.3788135593220339
$ws.Range("F3").Value = 0.3689149283290373
$ws.Range("G3").Value = 0.3788135593220339
$ws.Range("H3").Value = ""
